# Weekly roll: insert a new week's worth of data (3 rows: Primera/Segunda/
# Tercera) at the top of the data block (row 447), pushing all existing
# rows down by 3. The new rows largely mirror the previous top-of-block
# entries (same price figures), but carry a newer sample date and a
# slightly different "Tercera" volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 447; this shifts rows 447:515 down to
# 450:518 and extends the sheet dimension accordingly. Excel also copies
# row formatting (e.g. the date number format on column D) down from the
# row above onto the freshly inserted rows.
$ws.Rows("447:449").Insert()

# New date for the inserted week (2021-10-05, serial 44474).
$newDate = 44474

# --- Row 447: Primera ---
$ws.Cells.Item(447, 1).Value = 1
$ws.Cells.Item(447, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(447, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(447, 4).Value = $newDate
$ws.Cells.Item(447, 5).Value = 15
$ws.Cells.Item(447, 6).Value = 100112020
$ws.Cells.Item(447, 7).Value = "Tomate"
$ws.Cells.Item(447, 8).Value = "Larga vida"
$ws.Cells.Item(447, 9).Value = "Primera"
$ws.Cells.Item(447, 10).Value = 300
$ws.Cells.Item(447, 11).Value = 5000
$ws.Cells.Item(447, 12).Value = 5500
$ws.Cells.Item(447, 13).Value = 5250
$ws.Cells.Item(447, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(447, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(447, 16).Value = 525
$ws.Cells.Item(447, 17).Value = 10
$ws.Cells.Item(447, 18).Value = "Hortaliza"

# --- Row 448: Segunda ---
$ws.Cells.Item(448, 1).Value = 1
$ws.Cells.Item(448, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(448, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(448, 4).Value = $newDate
$ws.Cells.Item(448, 5).Value = 15
$ws.Cells.Item(448, 6).Value = 100112020
$ws.Cells.Item(448, 7).Value = "Tomate"
$ws.Cells.Item(448, 8).Value = "Larga vida"
$ws.Cells.Item(448, 9).Value = "Segunda"
$ws.Cells.Item(448, 10).Value = 350
$ws.Cells.Item(448, 11).Value = 4500
$ws.Cells.Item(448, 12).Value = 5000
$ws.Cells.Item(448, 13).Value = 4750
$ws.Cells.Item(448, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(448, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(448, 16).Value = 475
$ws.Cells.Item(448, 17).Value = 10
$ws.Cells.Item(448, 18).Value = "Hortaliza"

# --- Row 449: Tercera ---
$ws.Cells.Item(449, 1).Value = 1
$ws.Cells.Item(449, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(449, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(449, 4).Value = $newDate
$ws.Cells.Item(449, 5).Value = 15
$ws.Cells.Item(449, 6).Value = 100112020
$ws.Cells.Item(449, 7).Value = "Tomate"
$ws.Cells.Item(449, 8).Value = "Larga vida"
$ws.Cells.Item(449, 9).Value = "Tercera"
$ws.Cells.Item(449, 10).Value = 470
$ws.Cells.Item(449, 11).Value = 3500
$ws.Cells.Item(449, 12).Value = 4000
$ws.Cells.Item(449, 13).Value = 3750
$ws.Cells.Item(449, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(449, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(449, 16).Value = 375
$ws.Cells.Item(449, 17).Value = 10
$ws.Cells.Item(449, 18).Value = "Hortaliza"
